$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "26.257.53"
$ws.Cells.Item(2, 5).Value = "  +0.00%  "
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.666.20"
$ws.Cells.Item(3, 5).Value = "  +0.26%  "
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "1.008"
$ws.Cells.Item(4, 5).Value = "  -0.07%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "219.84"
$ws.Cells.Item(5, 5).Value = "  +0.56%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "0.5271"
$ws.Cells.Item(6, 5).Value = "  -0.76%  "
$ws.Cells.Item(7, 5).Value = "  -0.09%  "
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.2646"
$ws.Cells.Item(8, 5).Value = "  +0.31%  "
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.06354"
$ws.Cells.Item(9, 5).Value = "  -0.11%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "20.83"
$ws.Cells.Item(10, 5).Value = "  +1.33%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.07842"
$ws.Cells.Item(11, 5).Value = "  -0.01%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "4.499"
$ws.Cells.Item(12, 5).Value = "  -1.37%  "
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "1.663.16"
$ws.Cells.Item(13, 5).Value = "  -0.08%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "1.896.37"
$ws.Cells.Item(14, 5).Value = "  +0.36%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "0.5589"
$ws.Cells.Item(15, 5).Value = "  +0.71%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "0.0₅8091"
$ws.Cells.Item(16, 5).Value = "  -1.18%  "
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "65.68"
$ws.Cells.Item(17, 5).Value = "  -0.05%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "26.301.48"
$ws.Cells.Item(18, 5).Value = "  +0.05%  "
$ws.Cells.Item(19, 5).Value = "  -0.14%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "4.695"
$ws.Cells.Item(20, 5).Value = "  +0.38%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "198.66"
$ws.Cells.Item(21, 5).Value = "  +3.02%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "10.23"
$ws.Cells.Item(22, 5).Value = "  +0.18%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "6.034"
$ws.Cells.Item(23, 5).Value = "  -0.48%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "1.010"
$ws.Cells.Item(24, 5).Value = "  +0.00%  "
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "145.96"
$ws.Cells.Item(25, 5).Value = "  +0.84%  "
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "0.1215"
$ws.Cells.Item(26, 5).Value = "  -0.92%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "7.213"
$ws.Cells.Item(27, 5).Value = "  -0.34%  "
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "16.21"
$ws.Cells.Item(28, 5).Value = "  +0.28%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "1.526"
$ws.Cells.Item(29, 5).Value = "  +2.42%  "
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "0.05880"
$ws.Cells.Item(30, 5).Value = "  -0.85%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "1.285"
$ws.Cells.Item(31, 5).Value = "  +0.49%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "3.497"
$ws.Cells.Item(32, 5).Value = "  -2.60%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "3.328"
$ws.Cells.Item(33, 5).Value = "  +1.02%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "1.589"
$ws.Cells.Item(34, 5).Value = "  -1.95%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "0.9630"
$ws.Cells.Item(35, 5).Value = "  -0.11%  "
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "2.821"
$ws.Cells.Item(36, 5).Value = "  -0.05%  "
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "2.431"
$ws.Cells.Item(37, 5).Value = "  +0.28%  "
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.5794"
$ws.Cells.Item(38, 5).Value = "  +0.20%  "
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.01605"
$ws.Cells.Item(39, 5).Value = "  -0.15%  "
$ws.Cells.Item(40, 2).Value = "Maker"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "1.079.68"
$ws.Cells.Item(40, 5).Value = "  +2.70%  "
$ws.Cells.Item(41, 2).Value = "FraxShare"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "5.948"
$ws.Cells.Item(41, 5).Value = "  +0.46%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.8593"
$ws.Cells.Item(42, 5).Value = "  +0.30%  "
$ws.Cells.Item(43, 5).Value = "  -0.05%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "103.30"
$ws.Cells.Item(44, 5).Value = "  -0.89%  "
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "1.807.77"
$ws.Cells.Item(45, 5).Value = "  +0.30%  "
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "58.47"
$ws.Cells.Item(46, 5).Value = "  +1.91%  "
$ws.Cells.Item(47, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "0.0₈106"
$ws.Cells.Item(47, 5).Value = "  -0.61%  "
$ws.Cells.Item(48, 2).Value = "Frax"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "1.008"
$ws.Cells.Item(48, 5).Value = "  -1.24%  "
$ws.Cells.Item(49, 2).Value = "Mantle"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "0.4417"
$ws.Cells.Item(49, 5).Value = "  +0.98%  "
$ws.Cells.Item(50, 2).Value = "EnergySwap"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "8.040"
$ws.Cells.Item(50, 5).Value = "  +0.10%  "
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.05154"
$ws.Cells.Item(51, 5).Value = "  -0.19%  "
